$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.770599999999999
$ws.Range("A9").Value = -21.9514
$ws.Range("D12").Value = -6.932099999999997
$ws.Range("A13").Value = -22.0931
$ws.Range("D14").Value = -7.669300000000003
$ws.Range("A16").Value = -21.46379999999999
$ws.Range("A18").Value = -22.4045
$ws.Range("D19").Value = -8.0953
$ws.Range("A20").Value = -19.94339999999998
$ws.Range("A26").Value = -21.02609999999997
$ws.Range("D26").Value = -8.837999999999999
$ws.Range("A27").Value = -21.51819999999998
$ws.Range("D27").Value = -8.973599999999994
$ws.Range("A29").Value = -20.86339999999997
$ws.Range("D29").Value = -7.379900000000005
$ws.Range("A35").Value = -19.98519999999998
$ws.Range("A36").Value = -20.31999999999999
$ws.Range("D37").Value = -8.012099999999997
$ws.Range("D38").Value = -8.850099999999996
$ws.Range("A45").Value = -21.47789999999999
$ws.Range("D47").Value = -7.427900000000002
$ws.Range("D51").Value = -7.611499999999999
$ws.Range("D52").Value = -7.569899999999997
$ws.Range("A55").Value = -22.08509999999999
$ws.Range("D55").Value = -8.361200000000006
$ws.Range("A57").Value = -22.373
$ws.Range("A69").Value = -21.6475
$ws.Range("D69").Value = -7.209199999999997
$ws.Range("D70").Value = -7.300199999999998
$ws.Range("A76").Value = -19.91839999999998
$ws.Range("D76").Value = -7.716800000000001
$ws.Range("A78").Value = -19.96399999999998
$ws.Range("D81").Value = -7.423599999999999
$ws.Range("A82").Value = -22.07249999999999
$ws.Range("A83").Value = -22.03319999999999
$ws.Range("D83").Value = -8.385499999999997
$ws.Range("A93").Value = -21.00279999999998
$ws.Range("D94").Value = -6.799900000000001
$ws.Range("A97").Value = -22.0345
$ws.Range("D100").Value = -8.414299999999997
$ws.Range("D102").Value = -7.998599999999999
